# Update TrapCardData: refresh card effects for the new rules revision and
# fold the extra "交换机 / Trade machine" alt-card row into the main table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : 墓碑 / Tomb (unchanged name, effect text refreshed) -----------
$ws.Range("A2").Value = "墓碑"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "回合结束时：将主牌堆顶2张牌送墓。<br>`r`n开战时：用墓地顶端1张牌替换本牌。"
$ws.Range("E2").Value = "Tomb"

# --- Row 3 : 流沙 / Quicksand ----------------------------------------------
$ws.Range("A3").Value = "流沙"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "交锋时：本牌所在槽位和对位槽位的所有怪物牌点数变为1。"
$ws.Range("E3").Value = "Quicksand"
$ws.Range("F3").Value = "① Triggered (Before battle round): Check all rows above this card. If the total card count of a row is greater than the rank of this card, the ranks of all cards in that row are set to one."

# --- Row 4 : 落穴 / Pitfall ------------------------------------------------
$ws.Range("A4").Value = "落穴"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "交锋时：如果本牌所在槽位和对位槽位的怪物牌合计数量不小于本牌点数，则将那些怪物牌全部消灭，然后消灭本牌。"
$ws.Range("E4").Value = "Pitfall"
$ws.Range("F4").Value = "① Triggered (Pre battle): Check the row above this card. If there are more than one card in either column, send all cards in that row to Graveyard. Then send this card to Graveyard as well."

# --- Row 5 : 尖刺 / Spike trap (rank 1 -> 2) --------------------------------
$ws.Range("A5").Value = "尖刺"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "回合结束时：横置本牌，然后同槽位中所有其他牌点数减1。"
$ws.Range("E5").Value = "Spike trap"
$ws.Range("F5").Value = "① Triggered (Before battle round): Check all { Monster } cards above this card. Send those with the same rank as this card to Graveyard."

# --- Row 6 : 冷气喷口 / Cryogas vent (rank 2 -> 1) --------------------------
$ws.Range("A6").Value = "冷气喷口"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = "回合结束时：横置本牌所在槽位中所有牌。"
$ws.Range("E6").Value = "Cryogas vent"
$ws.Range("F6").Value = "① Lasting: Negate the effects of { Monster } cards in the row above this card."

# --- Row 7 : 吹箭 / Dart trap ----------------------------------------------
$ws.Range("A7").Value = "吹箭"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "有怪物牌移入本牌所在行或列时：消灭那张怪物牌，然后本牌点数减1。"
$ws.Range("E7").Value = "Dart trap"
$ws.Range("F7").Value = "① Triggered (Before battle round): Check the same row of this card. If there is any { Monster } card with lower rank than this card, select { Monster } cards with a total rank no higher than the rank of this card in this row. Send them to the Graveyard. Then send this card to the Graveyard as well."

# --- Row 8 : 滚石 / Boulder -------------------------------------------------
$ws.Range("A8").Value = "滚石"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = "死亡时：如果本牌点数大于1，则将本牌移动到相邻槽位中而不是送墓，然后使本牌和那个槽位中所有牌点数减1。<br>"
$ws.Range("E8").Value = "Boulder"
$ws.Range("F8").Value = "① Triggered(On sent to Graveyard from Battlefield): Select cards from below this card with a total rank no higher than this card. Send them to the Graveyard."

# --- Row 9 : 传送阵 / Teleporter --------------------------------------------
$ws.Range("A9").Value = "传送阵"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = "交锋时：如果本牌所在槽位和对位槽位的怪物牌合计数量大于1，则将那些怪物牌洗回主牌堆，然后消灭本牌。"
$ws.Range("E9").Value = "Teleporter"

# --- Row 10 : 爆桶 / Explosive barrel (rank 3 -> 2) -------------------------
$ws.Range("A10").Value = "爆桶"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = "点数降低时：本牌所在槽位和相邻槽位的所有牌点数减1，然后消灭本牌。"
$ws.Range("E10").Value = "Explosive barrel"

# --- Row 11 : 地雷 / Mine ----------------------------------------------------
$ws.Range("A11").Value = "地雷"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = "有牌移入本牌所在槽位时：该槽位中所有牌点数减1。"
$ws.Range("E11").Value = "Mine"
$ws.Range("F11").Value = "Triggered(Pre battle): Check all { monster } cards in the row above this card and in the same column with this card. If any of them has a rank higher than this card, send all the cards in that row to Graveyard. Then send this card to Graveyard as well."

# --- Row 12 : 石柱 / Pillar (rank 1 -> 2) -----------------------------------
$ws.Range("A12").Value = "石柱"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = "死亡时：选本牌所在行或列，其中的所有牌点数减1。"
$ws.Range("E12").Value = "Pillar"

# --- Row 13 : 延迟爆弹 / Delayed bomb (was the 暗门 row; rank 1 -> 2) -------
$ws.Range("A13").Value = "延迟爆弹"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = "回合结束时：移动到1个相邻槽位，然后点数减1，本牌点数因此降至0时，消灭本牌所在槽位的1张其他牌。"
$ws.Range("E13").Value = "Delayed bomb"

# --- Row 14 : 暗门 (was 交换机/Trade machine row; now folds in 暗门's data) -
$ws.Range("A14").Value = "暗门"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = "有怪物牌移入本牌所在槽位中时：将其移动到对位槽位。"
$ws.Range("E14").ClearContents()

# Row 15 (old 交换机 / Trade machine alt-card row) is dropped entirely.
$ws.Rows("15:15").Delete()

# Restore the view/selection state recorded for the refreshed sheet.
$ws.Range("D14").Select()
